# Apply cryptos list price/volume update (GitHub Actions data refresh)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Columns D cells below hold numeric-looking text (e.g. "27.00", "0.999") that
# must stay as TEXT (matching the sheet's existing inline-string cells) instead of
# being auto-coerced into numbers (which would drop formatting like trailing zeros).
# Pre-format them as Text before writing the values.
$textCells = @('D4','D5','D6','D8','D9','D10','D12','D14','D17','D19','D20','D21','D23','D24','D25','D26','D28','D31','D32','D33','D34','D35','D36','D37','D38','D40','D42','D43','D45','D49','D50','D51')
foreach ($addr in $textCells) {
    $ws.Range($addr).NumberFormat = "@"
}

# --- Row-by-row value updates ---

# Row 2
$ws.Range('D2').Value = '52.342.27'
$ws.Range('E2').Value = '  +5.68%  '
# Row 3
$ws.Range('D3').Value = '2.796.27'
$ws.Range('E3').Value = '  +5.99%  '
# Row 4
$ws.Range('D4').Value = '0.999'
$ws.Range('E4').Value = '  -0.06%  '
# Row 5
$ws.Range('D5').Value = '117.50'
$ws.Range('E5').Value = '  +5.01%  '
# Row 6
$ws.Range('D6').Value = '341.22'
$ws.Range('E6').Value = '  +4.69%  '
# Row 7
$ws.Range('E7').Value = '  +5.42%  '
# Row 8
$ws.Range('D8').Value = '0.999'
$ws.Range('E8').Value = '  -0.08%  '
# Row 9
$ws.Range('D9').Value = '0.581'
$ws.Range('E9').Value = '  +6.13%  '
# Row 10
$ws.Range('D10').Value = '42.28'
$ws.Range('E10').Value = '  +7.01%  '
# Row 11
$ws.Range('E11').Value = '  +7.54%  '
# Row 12
$ws.Range('D12').Value = '20.13'
$ws.Range('E12').Value = '  +0.19%  '
# Row 13
$ws.Range('E13').Value = '  +2.33%  '
# Row 14
$ws.Range('D14').Value = '7.64'
$ws.Range('E14').Value = '  +2.03%  '
# Row 15
$ws.Range('D15').Value = '3.232.64'
$ws.Range('E15').Value = '  +5.92%  '
# Row 16
$ws.Range('D16').Value = '2.792.27'
$ws.Range('E16').Value = '  +5.44%  '
# Row 17
$ws.Range('D17').Value = '0.889'
$ws.Range('E17').Value = '  +4.15%  '
# Row 18
$ws.Range('D18').Value = '52.112.59'
$ws.Range('E18').Value = '  +5.24%  '
# Row 19
$ws.Range('D19').Value = '3.25'
$ws.Range('E19').Value = '  +11.94%  '
# Row 20
$ws.Range('D20').Value = '13.42'
$ws.Range('E20').Value = '  +2.30%  '
# Row 21
$ws.Range('D21').Value = '6.97'
$ws.Range('E21').Value = '  +4.40%  '
# Row 22
$ws.Range('E22').Value = '  +4.28%  '
# Row 23
$ws.Range('D23').Value = '279.18'
$ws.Range('E23').Value = '  +3.98%  '
# Row 24
$ws.Range('D24').Value = '70.44'
$ws.Range('E24').Value = '  +1.96%  '
# Row 25
$ws.Range('D25').Value = '2.83'
$ws.Range('E25').Value = '  +10.26%  '
# Row 26
$ws.Range('D26').Value = '27.00'
$ws.Range('E26').Value = '  +3.59%  '
# Row 27
$ws.Range('E27').Value = '  -0.05%  '
# Row 28
$ws.Range('D28').Value = '10.23'
$ws.Range('E28').Value = '  +0.37%  '
# Row 29
$ws.Range('E29').Value = '  +1.16%  '
# Row 30
$ws.Range('E30').Value = '  +2.75%  '
# Row 31
$ws.Range('D31').Value = '35.08'
$ws.Range('E31').Value = '  +1.18%  '
# Row 32
$ws.Range('D32').Value = '50.40'
$ws.Range('E32').Value = '  +1.56%  '
# Row 33
$ws.Range('D33').Value = '5.73'
$ws.Range('E33').Value = '  +4.28%  '
# Row 34
$ws.Range('D34').Value = '0.0827'
$ws.Range('E34').Value = '  +2.37%  '
# Row 35
$ws.Range('D35').Value = '2.14'
$ws.Range('E35').Value = '  +5.12%  '
# Row 36
$ws.Range('D36').Value = '0.999'
$ws.Range('E36').Value = '  -0.14%  '
# Row 37
$ws.Range('D37').Value = '19.06'
$ws.Range('E37').Value = '  +0.23%  '
# Row 38
$ws.Range('D38').Value = '4.99'
$ws.Range('E38').Value = '  +0.58%  '
# Row 39
$ws.Range('E39').Value = '  +5.71%  '
# Row 40
$ws.Range('D40').Value = '2.76'
$ws.Range('E40').Value = '  +29.68%  '
# Row 41
$ws.Range('E41').Value = '  +12.50%  '
# Row 42
$ws.Range('B42').Value = 'EnergySwap'
$ws.Range('C42').Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range('D42').Value = '23.53'
$ws.Range('E42').Value = '  +3.58%  '
# Row 43
$ws.Range('B43').Value = 'WEMIXToken'
$ws.Range('C43').Value = 'https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix'
$ws.Range('D43').Value = '2.35'
$ws.Range('E43').Value = '  +5.34%  '
# Row 44
$ws.Range('E44').Value = '  +4.41%  '
# Row 45
$ws.Range('D45').Value = '127.26'
$ws.Range('E45').Value = '  -0.99%  '
# Row 46
$ws.Range('D46').Value = '2.110.59'
$ws.Range('E46').Value = '  +2.60%  '
# Row 47
$ws.Range('E47').Value = '  +2.72%  '
# Row 49
$ws.Range('D49').Value = '5.57'
$ws.Range('E49').Value = '  +6.98%  '
# Row 50
$ws.Range('D50').Value = '0.923'
$ws.Range('E50').Value = '  +22.86%  '
# Row 51
$ws.Range('D51').Value = '8.97'
$ws.Range('E51').Value = '  +1.12%  '

Write-Host "Updated 91 cells across the cryptos table"
